$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Scaled Feature Analysis"
